# Update cryptos list: Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.835.02"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "3.268.73"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'579.88"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "'185.19"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -4.46%  "
$ws.Range("D10").Value = "'6.58"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("D12").Value = "3.833.77"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'0.138"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "'27.48"
$ws.Range("E14").Value = "  -5.70%  "
$ws.Range("D15").Value = "67.898.02"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").Value = "3.230.21"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "'13.49"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "'396.58"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "'7.58"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'70.92"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E25").Value = "  -4.34%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'9.49"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").Value = "'22.62"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").Value = "'5.47"
$ws.Range("E31").Value = "  -5.71%  "
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -5.31%  "
$ws.Range("D35").Value = "'163.10"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").Value = "'27.08"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").Value = "'0.807"
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("D40").Value = "'4.51"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("D41").Value = "'6.30"
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").Value = "2.667.57"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "'40.72"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  -7.95%  "
$ws.Range("D46").Value = "'24.63"
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value = "'334.65"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").Value = "'0.972"
$ws.Range("E51").Value = "  -2.87%  "
